# Auto-generated Excel COM-interop edit script
# Applies the row-insertion + value-shift update described by the diff
# to sheets '展览' (index 1) and '全部类型' (index 4).

$wb = $excel.ActiveWorkbook
$targetSheets = @(1, 4)

foreach ($sheetIdx in $targetSheets) {
  $ws = $wb.Worksheets.Item($sheetIdx)

  # Row 2
  $ws.Range("B2").NumberFormat = "@"
  $ws.Range("B2").Value = "2024-03-23"
  $ws.Range("C2").Value = "苏州·Look Look动漫嘉年华"
  $ws.Range("D2").Value = "阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店"
  $ws.Range("E2").Value = "2024.03.23 10:00-03.23 17:30"
  $ws.Range("F2").Value = 841
  $ws.Range("G2").Value = 52.2
  $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81698"
  $ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202402/CP95X8ao1708934930351.jpeg"

  # Row 3
  $ws.Range("B3").NumberFormat = "@"
  $ws.Range("B3").Value = "2024-03-30"
  $ws.Range("C3").Value = "苏州·Anime Space动漫游戏展"
  $ws.Range("D3").Value = "金芳路与新发路交叉口东南120米 万龙大厦"
  $ws.Range("E3").Value = "2024.03.30 10:00-03.30 17:00"
  $ws.Range("F3").Value = 8
  $ws.Range("G3").Value = 55
  $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=82815"
  $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202403/XPHUZMJa1710327274324.jpeg"

  # Row 4
  $ws.Range("B4").NumberFormat = "@"
  $ws.Range("B4").Value = "2024-03-30"
  $ws.Range("C4").Value = "苏州·奇幻世界5.3动漫游戏展"
  $ws.Range("D4").Value = "龙河路1288号 乐动力苏州湾体育中心"
  $ws.Range("E4").Value = "2024.03.30 10:00-03.31 17:00"
  $ws.Range("F4").Value = 2157
  $ws.Range("G4").Value = 55
  $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=82002"
  $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg"

  # Row 5
  $ws.Range("B5").NumberFormat = "@"
  $ws.Range("B5").Value = "2024-03-31"
  $ws.Range("C5").Value = "张家港·META萌圆饿了"
  $ws.Range("D5").Value = "大新镇人民路18号 新香苑宴会厅"
  $ws.Range("E5").Value = "2024.03.31 10:00-03.31 17:00"
  $ws.Range("F5").Value = 53
  $ws.Range("G5").Value = 30
  $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=82407"
  $ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202403/pxZkSPcL1709707210340.jpeg"

  # Row 6
  $ws.Range("B6").NumberFormat = "@"
  $ws.Range("B6").Value = "2024-04-04"
  $ws.Range("C6").Value = "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会"
  $ws.Range("D6").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
  $ws.Range("E6").Value = "2024.04.04 10:00-04.05 17:00"
  $ws.Range("F6").Value = 12467
  $ws.Range("G6").Value = 60
  $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81827"
  $ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg"

  # Row 7
  $ws.Range("B7").NumberFormat = "@"
  $ws.Range("B7").Value = "2024-04-04"
  $ws.Range("C7").Value = "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会"
  $ws.Range("D7").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
  $ws.Range("E7").Value = "2024.04.04 10:00-04.05 17:00"
  $ws.Range("F7").Value = 12467
  $ws.Range("G7").Value = 60
  $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=81827"
  $ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg"

  # Row 8
  $ws.Range("B8").NumberFormat = "@"
  $ws.Range("B8").Value = "2024-04-04"
  $ws.Range("C8").Value = "常熟·SL动漫展"
  $ws.Range("D8").Value = "报慈北路218号 四季花园酒店(报慈北路店)"
  $ws.Range("E8").Value = "2024.04.04 10:00-04.04 17:00"
  $ws.Range("F8").Value = 56
  $ws.Range("G8").Value = 50
  $ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82250"
  $ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202403/rcvwgj7N1709518723294.jpeg"

  # Row 9
  $ws.Range("B9").NumberFormat = "@"
  $ws.Range("B9").Value = "2024-04-06"
  $ws.Range("C9").Value = "苏州·第一届寒假动漫展宅舞比赛-CF01"
  $ws.Range("D9").Value = "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店"
  $ws.Range("E9").Value = "2024.04.06 10:00-04.06 16:00"
  $ws.Range("F9").Value = 105
  $ws.Range("G9").Value = 49
  $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=80528"
  $ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"

  # Row 10
  $ws.Range("B10").NumberFormat = "@"
  $ws.Range("B10").Value = "2024-04-13"
  $ws.Range("C10").Value = "苏州·X-party 国漫游戏嘉年华03"
  $ws.Range("D10").Value = "秋枫街与开平路交叉口西南角 爱琴海购物中心"
  $ws.Range("E10").Value = "2024.04.13 10:00-04.14 17:00"
  $ws.Range("F10").Value = 500
  $ws.Range("G10").Value = 48
  $ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=82042"
  $ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg"

  # Row 11
  $ws.Range("B11").NumberFormat = "@"
  $ws.Range("B11").Value = "2024-04-13"
  $ws.Range("C11").Value = "苏州·绘时国乙1.0-秩序之外"
  $ws.Range("D11").Value = "兴中路与鲈乡北路交汇处 香漫商业广场"
  $ws.Range("E11").Value = "2024.04.13 13:00-04.13 20:00"
  $ws.Range("F11").Value = 451
  $ws.Range("G11").Value = 88
  $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=80789"
  $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202403/DI2ackIO1710137864319.jpeg"

  # Row 12
  $ws.Range("B12").NumberFormat = "@"
  $ws.Range("B12").Value = "2024-04-20"
  $ws.Range("C12").Value = "苏州·首届Redamancy动漫游戏嘉年华"
  $ws.Range("D12").Value = "清禾路886号 尹山湖大剧院"
  $ws.Range("E12").Value = "2024.04.20 10:00-04.20 17:00"
  $ws.Range("F12").Value = 1140
  $ws.Range("G12").Value = 60
  $ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=81879"
  $ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg"

  # Row 13
  $ws.Range("B13").NumberFormat = "@"
  $ws.Range("B13").Value = "2024-04-21"
  $ws.Range("C13").Value = "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0"
  $ws.Range("D13").Value = "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
  $ws.Range("E13").Value = "2024.04.21 10:00-04.21 21:00"
  $ws.Range("F13").Value = 931
  $ws.Range("G13").Value = 69.9
  $ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=78666"
  $ws.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202403/PlZCFPVs1710502485559.jpeg"

  # Row 14
  $ws.Range("B14").NumberFormat = "@"
  $ws.Range("B14").Value = "2024-05-01"
  $ws.Range("C14").Value = "昆山·第十二届理想乡动漫游戏展"
  $ws.Range("D14").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E14").Value = "2024.05.01 10:00-05.03 17:00"
  $ws.Range("F14").Value = 13632
  $ws.Range("G14").Value = 75
  $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=77196"
  $ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"

  # Row 15
  $ws.Range("B15").NumberFormat = "@"
  $ws.Range("B15").Value = "2024-05-01"
  $ws.Range("C15").Value = "苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区"
  $ws.Range("D15").Value = "金山南路288号 广电国际会展中心"
  $ws.Range("E15").Value = "2024.05.01 10:00-05.02 17:00"
  $ws.Range("F15").Value = 13894
  $ws.Range("G15").Value = 65
  $ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=79789"
  $ws.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"

  # Row 16
  $ws.Range("B16").NumberFormat = "@"
  $ws.Range("B16").Value = "2024-05-02"
  $ws.Range("C16").Value = "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会"
  $ws.Range("D16").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E16").Value = "2024.05.02 14:00-05.02 16:00"
  $ws.Range("F16").Value = 41
  $ws.Range("G16").Value = 1
  $ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=81116"
  $ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg"

  # Row 17
  $ws.Range("B17").NumberFormat = "@"
  $ws.Range("B17").Value = "2024-05-02"
  $ws.Range("C17").Value = "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会"
  $ws.Range("D17").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E17").Value = "2024.05.02 14:00-05.02 16:00"
  $ws.Range("F17").Value = 164
  $ws.Range("G17").Value = 1
  $ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=81100"
  $ws.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg"

  # Row 18
  $ws.Range("B18").NumberFormat = "@"
  $ws.Range("B18").Value = "2024-05-02"
  $ws.Range("C18").Value = "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会"
  $ws.Range("D18").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E18").Value = "2024.05.02 14:00-05.02 16:00"
  $ws.Range("F18").Value = 21
  $ws.Range("G18").Value = 1
  $ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=81119"
  $ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg"

  # Row 19
  $ws.Range("B19").NumberFormat = "@"
  $ws.Range("B19").Value = "2024-05-02"
  $ws.Range("C19").Value = "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会"
  $ws.Range("D19").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E19").Value = "2024.05.02 14:00-05.02 16:00"
  $ws.Range("F19").Value = 41
  $ws.Range("G19").Value = 1
  $ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=81118"
  $ws.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg"

  # Row 20
  $ws.Range("B20").NumberFormat = "@"
  $ws.Range("B20").Value = "2024-05-03"
  $ws.Range("C20").Value = "常熟·CDW·动漫展03"
  $ws.Range("D20").Value = "常熟国际展览中心 国际展览中心"
  $ws.Range("E20").Value = "2024.05.03 09:00-05.04 17:30"
  $ws.Range("F20").Value = 1041
  $ws.Range("G20").Value = 60
  $ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=82489"
  $ws.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202403/XK411blC1709794808211.jpeg"

  # Row 21
  $ws.Range("B21").NumberFormat = "@"
  $ws.Range("B21").Value = "2024-05-03"
  $ws.Range("C21").Value = "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会"
  $ws.Range("D21").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E21").Value = "2024.05.03 14:00-05.03 16:00"
  $ws.Range("F21").Value = 106
  $ws.Range("G21").Value = 1
  $ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=81120"
  $ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"

  # Row 22
  $ws.Range("B22").NumberFormat = "@"
  $ws.Range("B22").Value = "2024-05-03"
  $ws.Range("C22").Value = "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会"
  $ws.Range("D22").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
  $ws.Range("E22").Value = "2024.05.03 14:00-05.03 16:00"
  $ws.Range("F22").Value = 54
  $ws.Range("G22").Value = 1
  $ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=81114"
  $ws.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

  # Row 23
  $ws.Range("B23").NumberFormat = "@"
  $ws.Range("B23").Value = "2024-05-04"
  $ws.Range("C23").Value = "【大会员抢先购】苏州·OCG国潮动漫游戏嘉年华阿杰内场"
  $ws.Range("D23").Value = "苏州大道东688号 苏州国际博览中心"
  $ws.Range("E23").Value = "2024.05.04 09:00-05.04 17:00"
  $ws.Range("F23").Value = 391
  $ws.Range("G23").Value = 288
  $ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=82940"
  $ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202403/lLKmv48C1710511298160.jpeg"

  # Row 24
  $ws.Range("B24").NumberFormat = "@"
  $ws.Range("B24").Value = "2024-05-04"
  $ws.Range("C24").Value = "苏州·OCG国潮动漫游戏嘉年华"
  $ws.Range("D24").Value = "苏州大道东688号 苏州国际博览中心"
  $ws.Range("E24").Value = "2024.05.04 09:00-05.05 17:00"
  $ws.Range("F24").Value = 4990
  $ws.Range("G24").Value = 65
  $ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=82779"
  $ws.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"

  # Row 25
  $ws.Range("B25").NumberFormat = "@"
  $ws.Range("B25").Value = "2024-06-08"
  $ws.Range("C25").Value = "【会员购严选】苏州·Come in joy动漫国潮文化节"
  $ws.Range("D25").Value = "金山南路288号 广电国际会展中心"
  $ws.Range("E25").Value = "2024.06.08 10:00-06.09 17:00"
  $ws.Range("F25").Value = 237
  $ws.Range("G25").Value = 60
  $ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=82233"
  $ws.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg"

  # Row 25 is newly added: give A25 the same formatting as the existing index column (A2:A24)
  $ws.Range("A24").Copy() | Out-Null
  $ws.Range("A25").PasteSpecial(-4122) | Out-Null
  $ws.Range("A25").Value = 24

}
